# "Worked on front end" -----------------------------------------------------
# 1) Reorder the tabs so "listing an item" comes before "account page".
# 2) Tidy up a handful of TODO-list cells on the "account page" and
#    "main menu" sheets: a couple of split notes got merged into one note,
#    and one stray cell moved from column C to column E.
# 3) Leave the UI looking at the "main menu" tab afterwards.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) Move "listing an item" so it sits right before "account page" -----
$listingAnItem = $wb.Worksheets.Item("listing an item")
$accountPage   = $wb.Worksheets.Item("account page")
$listingAnItem.Move($accountPage)

# References above go stale once the sheet collection is reordered, so grab
# fresh handles to every sheet we still need to touch.
$accountPage   = $wb.Worksheets.Item("account page")
$listingAnItem = $wb.Worksheets.Item("listing an item")
$mainMenu      = $wb.Worksheets.Item("main menu")

# --- 2a) "account page" sheet ----------------------------------------------
# "button must remove item from the global and user objects" moves from C4
# to E4.
$accountPage.Range("E4").Value = $accountPage.Range("C4").Value2
$accountPage.Range("C4").Clear()

# "details button must lead to a page containing the details of the item
# being sold" note is dropped entirely.
$accountPage.Range("C5").Clear()

# The two purchase-history notes get merged into a single cell (E2); the
# old D2 cell disappears.
$accountPage.Range("E2").Value = "display personnalized purchase history //  display only users items currently for sale"
$accountPage.Range("D2").Clear()

# --- 2b) "main menu" sheet --------------------------------------------------
# The two search-bar notes get merged into a single cell (E5); D5 is
# cleared out (but keeps its formatting, unlike the cells above).
$mainMenu.Range("E5").Value = "clears the input in the search bar  // must also reset the search result // make button disabled if nothing was searched"
$mainMenu.Range("D5").ClearContents()

# --- 2c) selection on "listing an item" moves to C2 ------------------------
[void]$listingAnItem.Activate()
[void]$listingAnItem.Range("C2").Select()

# --- 3) Leave "main menu" as the active tab --------------------------------
$mainMenu.Activate()
